$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82: politeness_score (column B) becomes a real number instead of text
$ws.Cells.Item(82, 2).Value = 3

# Row 83: new annotation row appended after row 82
$ws.Cells.Item(83, 1).Value = "Ruilin"
$ws.Cells.Item(83, 2).NumberFormat = "@"
$ws.Cells.Item(83, 2).Value = "3"
$ws.Cells.Item(83, 3).Value = "无"
$ws.Cells.Item(83, 4).Value = "DIS"
$ws.Cells.Item(83, 5).Value = "RES"
$ws.Cells.Item(83, 6).Value = "21c11312-d736-4194-815f-bf7208ef5d55"
$ws.Cells.Item(83, 7).Value = "SJ60SbW0b_annotated.xlsx"
$ws.Cells.Item(83, 8).Value = "These sections include new experiments that illustrate the effect of varying the beta hyperparameter, demonstrate the strength of our approach on the larger scale Inception network for the ILSVRC 2014 classification challenge, and further highlight the effectiveness of our approach in diagnosing model failure modes."
